$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 2 = Jessica: "Work completed" cell was empty -> "[Not provided]"
#                   "Beans allocated" cell "100" -> "-"
$t.Cell(2, 2).Range.Text = "[Not provided]"
$t.Cell(2, 3).Range.Text = "-"

# Add a new "Unallocated points" row at the bottom of the table
$newRow = $t.Rows.Add()
$t.Cell($t.Rows.Count, 1).Range.Text = "Unallocated points"
$t.Cell($t.Rows.Count, 3).Range.Text = "100"
